$d = $word.ActiveDocument

# Hybrid bold + color (#2C3E50) highlighting for quantitative impact metrics.
# Word's Font.Color (an OLE_COLOR) is byte-order BGR, so RGB 2C3E50 -> BGR 503E2C.
$metricColor = 0x503E2C

# Map of (exact original paragraph text) -> (ordered list of metric substrings
# inside it that should become bold + colored). Matching on the paragraph's
# full original text keeps each lookup unambiguous even though some of the
# metric substrings (e.g. "23%", "71%") recur elsewhere in the document.
$targets = [ordered]@{
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" = @("23%", "64%")
    "• Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes" = @("±4.2%", "±2.1%", "71%", "87%")
    "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis" = @("73.5%", "`$4.7M")
    "• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion" = @("`$2")
    "• Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%" = @("±4.2%", "±2.1%")
    "• Increased voter turnout prediction accuracy from 71% to 87%" = @("71%", "87%")
    "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%" = @("34%", "28%")
}

foreach ($p in $d.Paragraphs) {
    $full = $p.Range.Text
    $full = $full.TrimEnd([char]13, [char]7)
    if ($targets.Contains($full)) {
        foreach ($metric in $targets[$full]) {
            $r = $p.Range.Duplicate
            $found = $r.Find.Execute($metric, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
            if ($found) {
                $r.Font.Bold = 1
                $r.Font.Color = $metricColor
            } else {
                Write-Output "MISS: $metric in $full"
            }
        }
    }
}
